$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.275.28'
$ws.Range("E2").Value = '  +5.89%  '

$ws.Range("D3").Value = '''3.284.50'
$ws.Range("E3").Value = '  +1.40%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = '''407.19'
$ws.Range("E5").Value = '  +3.06%  '

$ws.Range("D6").Value = '''111.84'
$ws.Range("E6").Value = '  +3.57%  '

$ws.Range("D7").Value = '''3.279.61'
$ws.Range("E7").Value = '  +1.47%  '

$ws.Range("D8").Value = '''0.566'
$ws.Range("E8").Value = '  -2.65%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").Value = '''0.618'
$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").Value = '''0.114'
$ws.Range("E11").Value = '  +14.67%  '

$ws.Range("D12").Value = '''38.62'
$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("D14").Value = '''3.799.64'
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").Value = '''8.15'
$ws.Range("E15").Value = '  -0.98%  '

$ws.Range("D16").Value = '''18.94'
$ws.Range("E16").Value = '  -0.50%  '

$ws.Range("D17").Value = '''3.337.36'
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").Value = '''60.103.72'
$ws.Range("E18").Value = '  +5.73%  '

$ws.Range("D19").Value = '''0.980'
$ws.Range("E19").Value = '  -4.74%  '

$ws.Range("D20").Value = '''10.59'
$ws.Range("E20").Value = '  -1.08%  '

$ws.Range("E21").Value = '  +3.56%  '

$ws.Range("D22").Value = '''3.29'
$ws.Range("E22").Value = '  -1.20%  '

$ws.Range("D23").Value = '''12.43'
$ws.Range("E23").Value = '  -3.14%  '

$ws.Range("D24").Value = '''295.57'
$ws.Range("E24").Value = '  +0.30%  '

$ws.Range("D25").Value = '''73.36'
$ws.Range("E25").Value = '  -0.98%  '

$ws.Range("D26").Value = '''3.08'
$ws.Range("E26").Value = '  -2.52%  '

$ws.Range("D27").Value = '''28.96'
$ws.Range("E27").Value = '  +3.42%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''7.38'
$ws.Range("E28").Value = '  +2.54%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.172'
$ws.Range("E29").Value = '  +2.57%  '

$ws.Range("B30").Value = 'LEO'
$ws.Range("C30").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D30").Value = '''4.22'
$ws.Range("E30").Value = '  -3.01%  '

$ws.Range("D31").Value = '''7.46'
$ws.Range("E31").Value = '  -1.34%  '

$ws.Range("E32").Value = '  +5.00%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("D34").Value = '''11.15'
$ws.Range("E34").Value = '  -0.79%  '

$ws.Range("D35").Value = '''2.45'
$ws.Range("E35").Value = '  +15.16%  '

$ws.Range("D36").Value = '''38.91'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").Value = '''0.0478'
$ws.Range("E37").Value = '  -0.38%  '

$ws.Range("D38").Value = '''52.10'
$ws.Range("E38").Value = '  +1.47%  '

$ws.Range("D39").Value = '''0.998'
$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("D40").Value = '''3.07'
$ws.Range("E40").Value = '  +6.66%  '

$ws.Range("D41").Value = '''3.30'
$ws.Range("E41").Value = '  -4.54%  '

$ws.Range("D42").Value = '''134.84'
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("E43").Value = '  +5.19%  '

$ws.Range("E44").Value = '  -2.20%  '

$ws.Range("E45").Value = '  -0.72%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''3.77'
$ws.Range("E46").Value = '  -4.08%  '

$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").Value = '''16.20'
$ws.Range("E47").Value = '  -4.24%  '

$ws.Range("E48").Value = '  +2.52%  '

$ws.Range("D49").Value = '''20.82'
$ws.Range("E49").Value = '  -5.47%  '

$ws.Range("D50").Value = '''2.118.25'
$ws.Range("E50").Value = '  -1.25%  '

$ws.Range("D51").Value = '''3.602.82'
$ws.Range("E51").Value = '  +0.88%  '
